$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cl = $nm.CustomLayouts.Item(1)
$tcs = $cl.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $c = $tcs.Item($i)
    Write-Host "Index $i : RGB=" $c.RGB
}
